$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# Row 2 (E2:O2) currently holds truly-empty text cells; the update fills
# them with the literal text "nan" (matching the rest of the sheet's
# "missing value" convention), keeping the cells as text, not numbers.
$ws.Range("E2:O2").Formula = '="nan"'
$ws.Range("E2:O2").Copy()
$ws.Range("E2:O2").PasteSpecial(-4163)

# Column A for rows 3-13 was mistakenly stored as "2" instead of "23"
# (this sheet is "Card23"); fix it back to "23" as text (not a number),
# matching the rest of the column/sheet's text-typed cells.
$ws.Range("A3:A13").Formula = '="23"'
$ws.Range("A3:A13").Copy()
$ws.Range("A3:A13").PasteSpecial(-4163)

$excel.CutCopyMode = 0
